$d = $word.ActiveDocument

# --- 1. Title block: collapse the 3 "Compact" paragraphs into a single
#        Heading1 paragraph "Robin F. De Guzman" wrapped by a bookmark.
$p3 = $d.Paragraphs(3)
$p3.Range.Delete()
$p1 = $d.Paragraphs(1)
$p1.Range.Delete()
$pTitle = $d.Paragraphs(1)
$pTitle.Style = "Heading1"
$titleStart = $pTitle.Range.Start
$titleEnd = $pTitle.Range.End - 1
$titleRange = $d.Range($titleStart, $titleEnd)
$d.Bookmarks.Add("robin-f.-de-guzman", $titleRange)

# --- 2. Every remaining Heading1/Heading2 bookmark currently wraps zero
#        characters (bookmarkStart immediately followed by bookmarkEnd,
#        *before* the run). Re-create each one spanning the heading's
#        text (excluding the trailing paragraph mark) so it wraps the
#        run, matching <bookmarkStart><r>...</r><bookmarkEnd>.
$bmNames = @()
foreach ($bm in $d.Bookmarks) {
    $bmNames += $bm.Name
}

foreach ($name in $bmNames) {
    $bm = $d.Bookmarks.Item($name)
    $para = $bm.Range.Paragraphs(1)
    $newStart = $para.Range.Start
    $newEnd = $para.Range.End - 1
    $newRange = $d.Range($newStart, $newEnd)
    $newName = $name
    if ($name -eq "st.-elizabeth-seton-school") {
        $newName = "st.elizabeth-seton-school"
    }
    $d.Bookmarks.Add($newName, $newRange)
}

# --- 3. Bold the three job-title lines (FirstParagraph runs).
$boldTargets = @(
    "Senior Software Engineer since July 2018",
    "Senior Developer from September 2012 to July 2018",
    "Technical Support Representative/Professional"
)
foreach ($t in $boldTargets) {
    $rng = $d.Content
    $rng.Find.Execute($t) | Out-Null
    $rng.Font.Bold = 1
}


# --- 4. Curly-quote / apostrophe normalisation in a few bullet runs.
$straightToCurly = @(
    @{ Old = "Provide tech support for client's customers"; New = "Provide tech support for client’s customers" },
    @{ Old = "Developed a system that accepts a csv file containing employee's time entries for a specified cutoff"; New = "Developed a system that accepts a csv file containing employee’s time entries for a specified cutoff" },
    @{ Old = "Developed a backend system that manages staff's team assignments, job positions/promotions, resignations. This system provides reset for forgotten passwords."; New = "Developed a backend system that manages staff’s team assignments, job positions/promotions, resignations. This system provides reset for forgotten passwords." }
)
foreach ($pair in $straightToCurly) {
    $rng = $d.Content
    $rng.Find.Execute($pair.Old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.New, 2) | Out-Null
}

Write-Output "done"
